# Apply the "Add files via upload" update to QE_holdings.xlsx
# 1) Bump the "as of" date in the confidential disclaimer text (A42)
#    from 2021-05-25 to 2021-05-26.
# 2) Update the Weight (col D) and Percent Change (col E) values for the
#    holdings rows (2-39) to the new model output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect it so the cells can be edited.
$ws.Unprotect()

# --- 1) Update the disclaimer text in A42 ---
$newText = "Model holdings provided as of 2021-05-26 for illustrative purposes only and are subject to change."
$ws.Range("A42").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`n" + $newText

# --- 2) Update Weight (D) / Percent Change (E) values ---
$ws.Range("D2").Value  = 0.07910549196550191
$ws.Range("E2").Value  = -0.0003940110323088497
$ws.Range("D3").Value  = 0.07298343114906863
$ws.Range("E3").Value  = -0.0009137136500873
$ws.Range("D4").Value  = 0.05828932543385014
$ws.Range("E4").Value  = 0.0004526462395544328
$ws.Range("D5").Value  = 0.04724627587922536
$ws.Range("E5").Value  = 0.001874779460272036
$ws.Range("D6").Value  = 0.04316118111848008
$ws.Range("E6").Value  = -0.008514540214828559
$ws.Range("D7").Value  = 0.04223395637346688
$ws.Range("E7").Value  = -0.0001235712079084017
$ws.Range("D8").Value  = 0.0394502525658147
$ws.Range("E8").Value  = -0.00593838193791163
$ws.Range("D9").Value  = 0.0330159275059858
$ws.Range("E9").Value  = -0.001194323450892321
$ws.Range("D10").Value = 0.03425440170808831
$ws.Range("E10").Value = 0.00738085463863869
$ws.Range("D11").Value = 0.03326368033419523
$ws.Range("E11").Value = -0.0003965953811890799
$ws.Range("D12").Value = 0.03166887116910921
$ws.Range("E12").Value = 0.0007141156867411791
$ws.Range("D13").Value = 0.02809319913412634
$ws.Range("E13").Value = 0.001702900607367885
$ws.Range("D14").Value = 0.03010218608165493
$ws.Range("E14").Value = -0.002239398974210882
$ws.Range("D15").Value = 0.02755144821603467
$ws.Range("E15").Value = 0.00183109707971596
$ws.Range("D16").Value = 0.02993982027264699
$ws.Range("E16").Value = 0.0003873623048056274
$ws.Range("D17").Value = 0.02661784481423902
$ws.Range("E17").Value = -0.01244485594466538
$ws.Range("D18").Value = 0.01978659333917464
$ws.Range("E18").Value = 0.0008791981712679142
$ws.Range("D19").Value = 0.02260943890435552
$ws.Range("E19").Value = 0.02064631956912022
$ws.Range("D20").Value = 0.02046968949278661
$ws.Range("E20").Value = -0.002124645892351285
$ws.Range("D21").Value = 0.02195939579043445
$ws.Range("E21").Value = 0.01167181599725375
$ws.Range("D22").Value = 0.0175323180444662
$ws.Range("E22").Value = 0.02388000463047169
$ws.Range("D23").Value = 0.02065148121108657
$ws.Range("E23").Value = 0.004380361379813902
$ws.Range("D24").Value = 0.01840184493949264
$ws.Range("E24").Value = 0.001016260162601812
$ws.Range("D25").Value = 0.01932080642458311
$ws.Range("E25").Value = -0.007503282686175172
$ws.Range("D26").Value = 0.01888082407606606
$ws.Range("E26").Value = -0.01393581081081086
$ws.Range("D27").Value = 0.01752506957084977
$ws.Range("E27").Value = 0.01913341274568192
$ws.Range("D28").Value = 0.01806957490891568
$ws.Range("E28").Value = 0.002406854722248886
$ws.Range("D29").Value = 0.01938560777871396
$ws.Range("E29").Value = -0.01042461225527591
$ws.Range("D30").Value = 0.0182847096058512
$ws.Range("E30").Value = 0.005454776100469516
$ws.Range("D31").Value = 0.01719917819705527
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 0.01826847302495041
$ws.Range("E32").Value = 0.001222067039106101
$ws.Range("D33").Value = 0.01685270115819011
$ws.Range("E33").Value = -0.003483870967741942
$ws.Range("D34").Value = 0.009073784242514212
$ws.Range("E34").Value = 0.003339138214759396
$ws.Range("D35").Value = 0.007267899525717875
$ws.Range("E35").Value = 0.002034547412933518
$ws.Range("D36").Value = 0.007499270803554187
$ws.Range("E36").Value = 0.01051614150396296
$ws.Range("D37").Value = 0.00732211810836874
$ws.Range("E37").Value = 0.003761780312029783
$ws.Range("D38").Value = 0.006661927131384674
$ws.Range("E38").Value = 0.002067284675980385
$ws.Range("E39").Value = 0.0005662507589152543

# Restore protection on the worksheet (it was protected before this edit).
$ws.Protect()
